# Add columns I (I0) and J (IF) to Sheet1, matching formatting of existing
# header column H, and fill in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
# Copy the formatting from H1 (bold, bordered, centered header style) onto
# the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (rows 2-17) ---
$iValues = @(7, 6, 6, 9, 7, 6, 4, 8, 7, 8, 7, 9, 9, 7, 5, 8)
$jValues = @(7, 6, 8, 9, 7, 6, 6, 8, 7, 8, 7, 9, 9, 8, 5, 8)

for ($idx = 0; $idx -lt 16; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
